# Proyecto_React_Carta_Gantt.xlsx - apply the tracked changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Row data: Fecha de Inicio (C), Fecha de Finalización (D) as real
# dates (short-date number format, centered) and Estado (E) as text
# with a status-coloured fill.
# ---------------------------------------------------------------

$rows = @(
    @{ R = 2;  C = 45495; D = 45495; E = "Completado" },
    @{ R = 3;  C = 45495; D = 45495; E = "Completado" },
    @{ R = 4;  C = 45495; D = 45495; E = "Completado" },
    @{ R = 5;  C = 45495; D = 45497; E = "Pendiente" },
    @{ R = 6;  C = 45495; D = 45496; E = "Completado" },
    @{ R = 7;  C = 45496; D = 45496; E = "pendiente" },
    @{ R = 8;  C = 45496; D = 45496; E = "Completado" },
    @{ R = 9;  C = 45496; D = 45496; E = "Completado" },
    @{ R = 10; C = 45496; D = 45497; E = "Pendiente" },
    @{ R = 11; C = 45496; D = 45497; E = "Completado" },
    @{ R = 12; C = 45496; D = 45497; E = "Completado" },
    @{ R = 13; C = 45496; D = 45497; E = "Pendiente" },
    @{ R = 14; C = 45496; D = 45497; E = "Pendiente" },
    @{ R = 15; C = 45496; D = 45497; E = "Por Finalizar" }
)

foreach ($row in $rows) {
    $r = $row.R

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.HorizontalAlignment = -4108
    $cCell.Value = $row.C
    $cCell.NumberFormat = "mm-dd-yy"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.HorizontalAlignment = -4108
    $dCell.Value = $row.D
    $dCell.NumberFormat = "mm-dd-yy"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $row.E
    if ($row.E -eq "Completado") {
        $eCell.Interior.ThemeColor = 7
    } elseif ($row.E -eq "Por Finalizar") {
        $eCell.Interior.Color = 255
    } else {
        $eCell.Interior.ThemeColor = 10
    }
}

# ---------------------------------------------------------------
# View tweaks: zoom + active selection
# ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 116
$ws.Range("H19").Select()
